$d = $word.ActiveDocument

# --- Edit 1: first paragraph gets two trailing spaces and a red parenthetical appended ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.SetRange($r1.Start, $r1.End - 1)  # exclude paragraph mark
$r1.Text = "This is a Microsoft word document.  "

# Insert the three red runs right after the text we just set, one w:r per
# sentence fragment (matches the original edit's run boundaries).
$ins = $d.Range($r1.End, $r1.End)
$ins.InsertAfter("(This is a change – Ve")
$ins.Font.Color = 255

$ins = $d.Range($ins.End, $ins.End)
$ins.InsertAfter("rsion for main branch")
$ins.Font.Color = 255

$ins = $d.Range($ins.End, $ins.End)
$ins.InsertAfter(")")
$ins.Font.Color = 255

# --- Edit 2: remove the trailing "ank God almighty, we are free at last." paragraph ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Delete()
